$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the "Late" column (old N, now becomes O).
$null = $ws.Columns("N:N").Insert()

# Excel copies the column width from the column to the left (M) onto the
# newly inserted column; M is width 11, so match that on the new N column.
$ws.Columns("N:N").ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab, and set its selection.
$ws.Activate()
$null = $ws.Range("N15").Select()
